$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.369.31"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.504.56"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "589.86"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "134.36"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.487"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "4.103.75"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "0.0000181"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "3.506.68"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "25.77"
$ws.Range("E16").Value = "  -5.98%  "
$ws.Range("D17").Value = "64.375.30"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "9.90"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "5.75"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "13.57"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").Value = "393.59"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").Value = "3.646.41"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "74.63"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "5.74"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  -7.55%  "
$ws.Range("D33").Value = "0.158"
$ws.Range("E33").Value = "  +8.20%  "
$ws.Range("D34").Value = "3.530.25"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "23.39"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "6.95"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").Value = "167.64"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").Value = "0.0789"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.44"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "24.88"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").Value = "1.67"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("D47").Value = "1.17"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "2.384.19"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50").Value = "0.898"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -0.71%  "
